$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 152, shifting existing rows 152:184 down to 153:185.
$ws.Rows.Item(152).Insert()

# Populate the new row 152 with the new weekly record (same market/category, new week).
$ws.Cells.Item(152, 1).Value = 7
$ws.Cells.Item(152, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(152, 3).Value = "Ñuble"
$ws.Cells.Item(152, 4).Value = 44637
$ws.Cells.Item(152, 5).Value = 16
$ws.Cells.Item(152, 6).Value = 100112017
$ws.Cells.Item(152, 7).Value = "Apio"
$ws.Cells.Item(152, 8).Value = "Americana (o)"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 60
$ws.Cells.Item(152, 11).Value = 8000
$ws.Cells.Item(152, 12).Value = 8500
$ws.Cells.Item(152, 13).Value = 8250
$ws.Cells.Item(152, 14).Value = "$/docena de matas"
$ws.Cells.Item(152, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(152, 16).Value = 1375
$ws.Cells.Item(152, 17).Value = 6
$ws.Cells.Item(152, 18).Value = "Hortaliza"
